$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# For Price values that look like plain numbers, force the cell to stay
# text (matching the original inlineStr cells) by temporarily applying a
# text number format, then clearing the format again so no extra styling
# is left behind on the cell.

$ws.Range("D2").Value = "67.154.12"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "2.470.92"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.51"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.29%  "

$ws.Range("E8").Value = "  -0.48%  "

$ws.Range("E9").Value = "  +2.08%  "

$ws.Range("E10").Value = "  +0.22%  "

$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("E12").Value = "  +0.97%  "

$ws.Range("D13").Value = "2.917.84"
$ws.Range("E13").Value = "  +0.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.37"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.80%  "

$ws.Range("D15").Value = "67.114.49"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("D17").Value = "2.428.37"
$ws.Range("E17").Value = "  -0.64%  "

$ws.Range("E18").Value = "  -1.77%  "

$ws.Range("E19").Value = "  -1.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "348.40"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.98"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.44"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("E24").Value = "  -1.59%  "

$ws.Range("E25").Value = "  -0.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.21"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.34%  "

$ws.Range("D27").Value = "2.595.27"
$ws.Range("E27").Value = "  +0.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("D29").Value = "0.0₃0899"
$ws.Range("E29").Value = "  -0.61%  "

$ws.Range("E30").Value = "  -3.06%  "

$ws.Range("E31").Value = "  -0.25%  "

$ws.Range("E32").Value = "  -0.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.76"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.03%  "

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("E35").Value = "  +2.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.80"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.85%  "

$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("E38").Value = "  -1.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.32"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.68%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.81"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.36%  "

$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "142.75"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.49%  "

$ws.Range("E46").Value = "  +0.60%  "

$ws.Range("E47").Value = "  +0.18%  "

$ws.Range("E48").Value = "  -1.11%  "

$ws.Range("E49").Value = "  +1.39%  "

$ws.Range("E50").Value = "  -1.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.582"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.01%  "
